# Week 9 second-half agenda update
# Rearranges the deck: old "Second Half" overview becomes slide 1 with the
# updated week list, old slide 1 (title slide) is replaced by a "This Week"
# agenda slide, and four new status slides are appended.

$p = $ppt.ActivePresentation
$layoutTitle   = $p.SlideMaster.CustomLayouts.Item(1)   ; # Title Slide
$layoutContent = $p.SlideMaster.CustomLayouts.Item(2)   ; # Title and Content

# ---------------------------------------------------------------------
# Slide 1: was an empty title slide -> becomes the "Second Half" agenda
# ---------------------------------------------------------------------
$p.Slides.Item(1).Delete()
$s1 = $p.Slides.AddSlide(1, $layoutContent)

$s1.Shapes.Item(1).TextFrame.TextRange.Text = "Second Half"

$tr1 = $s1.Shapes.Item(2).TextFrame.TextRange
$tr1.Text = "Week 9 (3/21) – "
$cur = $tr1.InsertAfter("Github")
$cur = $cur.InsertAfter(" & Raster and Imagery Data")
$cur = $cur.InsertAfter("`rWeek 10 (3/28) – Imagery Data Management, Lidar, and Mosaics")
$cur = $cur.InsertAfter("`rWeek 11 (4/4) – ")
$cur = $cur.InsertAfter("Arcpy.mapping")
$cur = $cur.InsertAfter(", working with MXDs")
$cur = $cur.InsertAfter("`rWeek 12 (4/11) – Network and Spatial Analysis & Demo of Pandas ")
$cur = $cur.InsertAfter("`rWeek 13 (4/18) – ")
$cur = $cur.InsertAfter("Jupyter")
$cur = $cur.InsertAfter(" and Projects in class")
$cur = $cur.InsertAfter("`rWeek 14 (4/25) – HTML/")
$cur = $cur.InsertAfter("Javascript")
$cur = $cur.InsertAfter(" and web development")
$cur = $cur.InsertAfter("`rWeek 15 (5/2) – Final Project Presentations")
$cur = $cur.InsertAfter("`rFinal Period (5/9) – Final Project Presentations")

# ---------------------------------------------------------------------
# Slide 2: was "Second Half" overview -> becomes "This Week" agenda
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "This Week"

$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange
$tr2.Text = "Install Source Tree"
$cur = $tr2.InsertAfter("`r`rWalk through GitHub & ")
$cur = $cur.InsertAfter("SourceTree")
$cur = $cur.InsertAfter(" & Slack together")
$cur = $cur.InsertAfter("`r`rRaster Lecture")
$cur = $cur.InsertAfter("`r`rRaster Exercise")
$cur = $cur.InsertAfter("`r")

# ---------------------------------------------------------------------
# Slide 3 (new): Next Week
# ---------------------------------------------------------------------
$s3 = $p.Slides.AddSlide(3, $layoutContent)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Next Week"
$s3.Shapes.Item(2).TextFrame.TextRange.Text = "Imagery, Mosaic Datasets, and LAS (Lidar) Datasets!"

# ---------------------------------------------------------------------
# Slide 4 (new): Project 2 is available!
# ---------------------------------------------------------------------
$s4 = $p.Slides.AddSlide(4, $layoutContent)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Project 2 is available!"
$tr4 = $s4.Shapes.Item(2).TextFrame.TextRange
$tr4.Text = "Go to webpage"
$cur = $tr4.InsertAfter("`r")

# ---------------------------------------------------------------------
# Slide 5 (new): Final Project is Available
# ---------------------------------------------------------------------
$s5 = $p.Slides.AddSlide(5, $layoutContent)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Final Project is Available"
$url5 = "https://blackboard.slu.edu/webapps/blackboard/content/listContentEditable.jsp?content_id=_2433173_1&course_id=_153733_1&mode=reset"
$tr5 = $s5.Shapes.Item(2).TextFrame.TextRange
$tr5.Text = $url5
$cur = $tr5.InsertAfter("`r`rWould like to see you use GitHub for submission though!")
$hl5 = $s5.Shapes.Item(2).TextFrame.TextRange.Characters(1, $url5.Length)
$hl5.ActionSettings.Item(1).Hyperlink.Address = $url5

# ---------------------------------------------------------------------
# Slide 6 (new): Week 9 Exercise 9
# ---------------------------------------------------------------------
$s6 = $p.Slides.AddSlide(6, $layoutContent)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "Week 9,  `"Exercise 9`""
$url6 = "https://blackboard.slu.edu/webapps/assignment/uploadAssignment?content_id=_2483732_1&course_id=_153733_1&assign_group_id=&mode=cpview"
$tr6 = $s6.Shapes.Item(2).TextFrame.TextRange
$tr6.Text = $url6
$cur = $tr6.InsertAfter("`r")
$hl6 = $s6.Shapes.Item(2).TextFrame.TextRange.Characters(1, $url6.Length)
$hl6.ActionSettings.Item(1).Hyperlink.Address = $url6

Write-Host "Final slide count:" $p.Slides.Count
